$p = $ppt.ActivePresentation

# The license/credits slide is the second slide in the deck.
$s = $p.Slides.Item(2)

# Locate the "Title 1" placeholder that holds the license paragraph.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Title 1") {
        $sh = $candidate
        break
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(1)
}

$tr = $sh.TextFrame.TextRange

# "CC BY-NC 4.0. To view a copy of this license, visit " -> "CC BY-SA 4.0. ..."
$t = $tr.Text
$idx = $t.IndexOf("BY-NC ")
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, 6)
    $sub.Text = "BY-SA "
}

# Hyperlink display text: ".../licenses/by-nc/4.0" -> ".../licenses/by-sa/4.0"
$t = $tr.Text
$idxUrl = $t.IndexOf("creativecommons.org/licenses/by-nc/4.0")
if ($idxUrl -ge 0) {
    $subUrl = $tr.Characters($idxUrl + 1, 39)
    $subUrl.Text = "creativecommons.org/licenses/by-sa/4.0"
}
